# Apply the edits described by the diff:
#  1. Table preferred width: auto -> 100% (pct)
#  2. Twelve numeric cell value updates in the coefficient table

$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# --- 1. Table width: auto/0 -> pct 100% -------------------------------
# wdPreferredWidthPercent = 2. Setting PreferredWidth = 250 with a
# Percent width type yields a stored width of 250*20 = 5000
# (fiftieths-of-a-percent), i.e. 100%, matching the target change from
# <w:tblW w:type="auto" w:w="0"/> to a 100% "pct" width.
$tbl.PreferredWidthType = 2
$tbl.PreferredWidth = 250

# --- 2. Numeric cell replacements --------------------------------------
$replacements = @(
    @("-0.1842", "-0.1337"),
    @("0.1435",  "0.0485"),
    @("1.2842",  "2.7569"),
    @("0.1991",  "0.0058"),
    @("0.0069",  "0.0047"),
    @("0.0062",  "0.0045"),
    @("1.1237",  "1.0406"),
    @("0.2611",  "0.2981"),
    @("-0.0016", "-0.0002"),
    @("0.0043",  "0.0010"),
    @("0.3810",  "0.2492"),
    @("0.7032",  "0.8032")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
